# Fruta / hortaliza, semanal
# Insert 3 new weekly rows of data right before the existing row 447,
# shifting all subsequent rows down by 3 (old row N -> new row N+3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 447 (pushes existing 447.. down to 450..)
$ws.Range("A447:A449").EntireRow.Insert()

# --- New row 447 --------------------------------------------------------
$ws.Cells.Item(447, 1).Value  = 8
$ws.Cells.Item(447, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(447, 3).Value  = "Coquimbo"
$ws.Cells.Item(447, 4).Value  = 44641
$ws.Cells.Item(447, 5).Value  = 4
$ws.Cells.Item(447, 6).Value  = "Fruta"
$ws.Cells.Item(447, 7).Value  = 100104
$ws.Cells.Item(447, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(447, 9).Value  = 100104005
$ws.Cells.Item(447, 10).Value = "Pera"
$ws.Cells.Item(447, 11).Value = "Packham's Triumph"
$ws.Cells.Item(447, 12).Value = "Especial"
$ws.Cells.Item(447, 13).Value = 16
$ws.Cells.Item(447, 14).Value = 235000
$ws.Cells.Item(447, 15).Value = 240000
$ws.Cells.Item(447, 16).Value = 237500
$ws.Cells.Item(447, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(447, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(447, 19).Value = 528
$ws.Cells.Item(447, 20).Value = 450

# --- New row 448 --------------------------------------------------------
$ws.Cells.Item(448, 1).Value  = 8
$ws.Cells.Item(448, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(448, 3).Value  = "Coquimbo"
$ws.Cells.Item(448, 4).Value  = 44641
$ws.Cells.Item(448, 5).Value  = 4
$ws.Cells.Item(448, 6).Value  = "Fruta"
$ws.Cells.Item(448, 7).Value  = 100104
$ws.Cells.Item(448, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(448, 9).Value  = 100104005
$ws.Cells.Item(448, 10).Value = "Pera"
$ws.Cells.Item(448, 11).Value = "Packham's Triumph"
$ws.Cells.Item(448, 12).Value = "Primera"
$ws.Cells.Item(448, 13).Value = 20
$ws.Cells.Item(448, 14).Value = 205000
$ws.Cells.Item(448, 15).Value = 210000
$ws.Cells.Item(448, 16).Value = 207500
$ws.Cells.Item(448, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(448, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(448, 19).Value = 461
$ws.Cells.Item(448, 20).Value = 450

# --- New row 449 --------------------------------------------------------
$ws.Cells.Item(449, 1).Value  = 8
$ws.Cells.Item(449, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(449, 3).Value  = "Coquimbo"
$ws.Cells.Item(449, 4).Value  = 44641
$ws.Cells.Item(449, 5).Value  = 4
$ws.Cells.Item(449, 6).Value  = "Fruta"
$ws.Cells.Item(449, 7).Value  = 100104
$ws.Cells.Item(449, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(449, 9).Value  = 100104005
$ws.Cells.Item(449, 10).Value = "Pera"
$ws.Cells.Item(449, 11).Value = "Packham's Triumph"
$ws.Cells.Item(449, 12).Value = "Segunda"
$ws.Cells.Item(449, 13).Value = 16
$ws.Cells.Item(449, 14).Value = 175000
$ws.Cells.Item(449, 15).Value = 180000
$ws.Cells.Item(449, 16).Value = 177500
$ws.Cells.Item(449, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(449, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(449, 19).Value = 394
$ws.Cells.Item(449, 20).Value = 450
